# Optuna Attempt (go back with original)
# Update forecast metrics on "Forecast Comparison" and roll-up totals on "Summary".

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet ---

# Row 2 (W8)
$wsForecast.Range("D2").Value = 23
$wsForecast.Range("H2").Value = 16.35
$wsForecast.Range("L2").Value = 0.83

# Row 3 (W9)
$wsForecast.Range("D3").Value = 25
$wsForecast.Range("H3").Value = 14.12
$wsForecast.Range("L3").Value = 0.83

# Row 4 (W10)
$wsForecast.Range("H4").Value = 12.01
$wsForecast.Range("L4").Value = 0.84

# Row 5 (W11)
$wsForecast.Range("H5").Value = 9.640000000000001
$wsForecast.Range("L5").Value = 0.8100000000000001

# Row 6 (W12)
$wsForecast.Range("H6").Value = 9.01
$wsForecast.Range("L6").Value = 0.84

# Row 7 (W13)
$wsForecast.Range("H7").Value = 8.01
$wsForecast.Range("L7").Value = 1.05

# Row 8 (W14)
$wsForecast.Range("H8").Value = 6.72
$wsForecast.Range("L8").Value = 0.9399999999999999

# Row 9 (W15)
$wsForecast.Range("H9").Value = 5.49
$wsForecast.Range("L9").Value = 1.01

# Row 10 (W16)
$wsForecast.Range("H10").Value = 4.88
$wsForecast.Range("L10").Value = 1.05

# Row 11 (W17)
$wsForecast.Range("H11").Value = 3.57
$wsForecast.Range("L11").Value = 1.02

# Row 12 (W18)
$wsForecast.Range("H12").Value = 2.8
$wsForecast.Range("L12").Value = 0.82

# Row 13 (W19)
$wsForecast.Range("H13").Value = 1.72

# Row 14 (W20)
$wsForecast.Range("H14").Value = 0.75
$wsForecast.Range("I14").Value = "Low"
$wsForecast.Range("L14").Value = 1.05

# Row 15 (W21)
$wsForecast.Range("L15").Value = 1.16

# Row 16 (W22)
$wsForecast.Range("L16").Value = 1.06

# Row 17 (W23)
$wsForecast.Range("L17").Value = 1.13

# --- Summary sheet ---
# These cells hold plain-text numbers (not numeric values), so force the
# Text number format before assigning to avoid COM's automatic
# string->number coercion.

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "472"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "230"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "106"

$wsSummary.Range("B12").NumberFormat = "@"
$wsSummary.Range("B12").Value = "32"

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "23"
